$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 16 ("Fiscal Impact" / current) - update H16:T16
$row16 = @{
    "H16" = -2.4186
    "I16" = -0.9767
    "J16" = 0.0753
    "K16" = -0.3713
    "L16" = -0.0351
    "M16" = 0.2572
    "N16" = 0.0139
    "O16" = 0.9707
    "P16" = 0.0487
    "Q16" = 0.0594
    "R16" = -0.323
    "S16" = -0.6118
    "T16" = 0.072
}

foreach ($addr in $row16.Keys) {
    $ws.Range($addr).Value = $row16[$addr]
}

# Row 44 ("Fiscal Impact" / difference) - update H44:T44
$row44 = @{
    "H44" = 0.0862
    "I44" = -0.077
    "J44" = 0.0826
    "K44" = 0.369
    "L44" = 0.2694
    "M44" = 0.5498
    "N44" = 0.5617
    "O44" = 0.9883
    "P44" = 0.1151
    "Q44" = 0.0016
    "R44" = -0.2508
    "S44" = 72.098
    "T44" = 1.9566
}

foreach ($addr in $row44.Keys) {
    $ws.Range($addr).Value = $row44[$addr]
}
